$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds text-formatted numeric-looking strings (e.g. '48.237.58',
# '1.00') in the source data. Pre-set the whole column to Text format so Excel's
# COM layer doesn't auto-convert these into real numbers on assignment (which would
# strip formatting like trailing zeros / thousand-separator dots).
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "48.237.58"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3
$ws.Range("D3").Value = "2.506.75"
$ws.Range("E3").Value = "  +0.03%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").Value = "319.54"
$ws.Range("E5").Value = "  -0.91%  "

# Row 6
$ws.Range("D6").Value = "106.89"
$ws.Range("E6").Value = "  -0.92%  "

# Row 7
$ws.Range("D7").Value = "0.524"
$ws.Range("E7").Value = "  +0.19%  "

# Row 8
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -3.95%  "

# Row 10
$ws.Range("D10").Value = "38.99"
$ws.Range("E10").Value = "  -2.97%  "

# Row 11
$ws.Range("D11").Value = "19.85"
$ws.Range("E11").Value = "  +1.41%  "

# Row 12
$ws.Range("D12").Value = "0.0806"
$ws.Range("E12").Value = "  -0.99%  "

# Row 13
$ws.Range("E13").Value = "  -0.55%  "

# Row 14
$ws.Range("D14").Value = "7.07"
$ws.Range("E14").Value = "  -1.40%  "

# Row 15
$ws.Range("D15").Value = "2.902.55"
$ws.Range("E15").Value = "  +0.15%  "

# Row 16
$ws.Range("D16").Value = "2.510.71"
$ws.Range("E16").Value = "  +0.19%  "

# Row 17
$ws.Range("D17").Value = "0.832"
$ws.Range("E17").Value = "  -2.12%  "

# Row 18
$ws.Range("D18").Value = "48.143.62"

# Row 19
$ws.Range("D19").Value = "12.96"
$ws.Range("E19").Value = "  -2.78%  "

# Row 20
$ws.Range("D20").Value = "2.94"
$ws.Range("E20").Value = "  +6.20%  "

# Row 21
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  +0.43%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0937"
$ws.Range("E22").Value = "  -0.35%  "

# Row 23
$ws.Range("D23").Value = "71.13"
$ws.Range("E23").Value = "  +0.32%  "

# Row 24
$ws.Range("D24").Value = "272.83"
$ws.Range("E24").Value = "  +10.36%  "

# Row 25
$ws.Range("E25").Value = "  -2.10%  "

# Row 26
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("D27").Value = "25.85"
$ws.Range("E27").Value = "  +0.47%  "

# Row 28
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.144"
$ws.Range("E28").Value = "  +1.77%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.20"
$ws.Range("E29").Value = "  +5.87%  "

# Row 30
$ws.Range("D30").Value = "9.72"
$ws.Range("E30").Value = "  -4.89%  "

# Row 31
$ws.Range("D31").Value = "34.58"
$ws.Range("E31").Value = "  -0.61%  "

# Row 32
$ws.Range("D32").Value = "49.47"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").Value = "19.18"
$ws.Range("E33").Value = "  -4.19%  "

# Row 34
$ws.Range("E34").Value = "  -0.16%  "

# Row 35
$ws.Range("D35").Value = "5.29"
$ws.Range("E35").Value = "  -1.68%  "

# Row 36
$ws.Range("D36").Value = "0.0778"
$ws.Range("E36").Value = "  -0.58%  "

# Row 37
$ws.Range("E37").Value = "  -0.93%  "

# Row 38
$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  -1.62%  "

# Row 39
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").Value = "  -3.43%  "

# Row 40
$ws.Range("E40").Value = "  -1.28%  "

# Row 41
$ws.Range("E41").Value = "  +1.23%  "

# Row 42
$ws.Range("D42").Value = "120.26"

# Row 43
$ws.Range("D43").Value = "21.79"
$ws.Range("E43").Value = "  -2.02%  "

# Row 44
$ws.Range("E44").Value = "  +2.35%  "

# Row 45
$ws.Range("D45").Value = "2.000.82"
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
$ws.Range("D46").Value = "3.18"
$ws.Range("E46").Value = "  +2.98%  "

# Row 47
$ws.Range("E47").Value = "  +5.07%  "

# Row 48
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").Value = "  -1.07%  "

# Row 49
$ws.Range("D49").Value = "8.94"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  +1.32%  "

# Row 51
$ws.Range("D51").Value = "78.85"
$ws.Range("E51").Value = "  +1.94%  "
